# Weekly update: insert a new week of "Coliflor" (cauliflower) price data
# at the top of the Terminal La Palmera de La Serena block (rows 293-294),
# pushing the existing rows 293:392 down to 295:394.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 293 - this shifts all existing data
# (old rows 293:392) down to 295:394, matching the rest of the dataset.
$ws.Rows("293:294").Insert()

# --- New row 293: "Primera" quality ---
$ws.Cells.Item(293, 1).Value = 8
$ws.Cells.Item(293, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(293, 3).Value = "Coquimbo"
$ws.Cells.Item(293, 4).Value = 44468
$ws.Cells.Item(293, 5).Value = 4
$ws.Cells.Item(293, 6).Value = 100112008
$ws.Cells.Item(293, 7).Value = "Coliflor"
$ws.Cells.Item(293, 8).Value = "Sin especificar"
$ws.Cells.Item(293, 9).Value = "Primera"
$ws.Cells.Item(293, 10).Value = 3200
$ws.Cells.Item(293, 11).Value = 600
$ws.Cells.Item(293, 12).Value = 700
$ws.Cells.Item(293, 13).Value = 650
$ws.Cells.Item(293, 14).Value = "$/unidad"
$ws.Cells.Item(293, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(293, 16).Value = 650
$ws.Cells.Item(293, 17).Value = 1
$ws.Cells.Item(293, 18).Value = "Hortaliza"

# --- New row 294: "Segunda" quality ---
$ws.Cells.Item(294, 1).Value = 8
$ws.Cells.Item(294, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(294, 3).Value = "Coquimbo"
$ws.Cells.Item(294, 4).Value = 44468
$ws.Cells.Item(294, 5).Value = 4
$ws.Cells.Item(294, 6).Value = 100112008
$ws.Cells.Item(294, 7).Value = "Coliflor"
$ws.Cells.Item(294, 8).Value = "Sin especificar"
$ws.Cells.Item(294, 9).Value = "Segunda"
$ws.Cells.Item(294, 10).Value = 1600
$ws.Cells.Item(294, 11).Value = 500
$ws.Cells.Item(294, 12).Value = 550
$ws.Cells.Item(294, 13).Value = 525
$ws.Cells.Item(294, 14).Value = "$/unidad"
$ws.Cells.Item(294, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(294, 16).Value = 525
$ws.Cells.Item(294, 17).Value = 1
$ws.Cells.Item(294, 18).Value = "Hortaliza"
